# Apply the "POPRAWIONY MANGER, FOLDER DELEM" changes to Sheet1 / ListaElementow

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks -------------------------------------------------
# The stored (OOXML) column width and the COM ColumnWidth property differ by
# a constant padding offset (~0.8333 chars) for this runtime/font, so the
# COM-side values below are chosen to land exactly on width=5 / width=25
# once persisted.
$ws.Columns.Item(1).ColumnWidth = 4.166666666666667
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668

# --- Remove the two trailing rows (13 & 14) so the table ends at row 12 --
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()

# --- New data for rows 2-12 (A = ID, C = Nazwa pliku, D = Status) --------
$data = @(
    @(2,  100, "man.manualmode.dld",      "W trakcie"),
    @(3,  100, "prd.2_600_250_20_15.dld", "Archiwum"),
    @(4,  100, "prd.3g_4_100_50_50.dld",  "W trakcie"),
    @(5,  100, "prd.40034102p4p5.dld",    "Archiwum"),
    @(6,  100, "prd.40662901siatka.dld",  "Archiwum"),
    @(7,  100, "prd.4_100k9050.dld",      "Archiwum"),
    @(8,  100, "prd.8_300_100_50_.dld",   "Archiwum"),
    @(9,  100, "prd.GRAFIKA.dld",         "Archiwum"),
    @(10, 100, "prd.testTrakcie.dld",     "Archiwum"),
    @(11, 101, "prd.TEST1.dld",           "W trakcie"),
    @(12, 102, "prd.626246TESAT.dld",     "Gotowe")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
